$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.866.22"
$ws.Range("E2").Value = "  -1.05%  "

$ws.Range("D3").Value = "3.266.06"
$ws.Range("E3").Value = "  +0.05%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.17"
$ws.Range("E5").Value = "  -0.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.24"
$ws.Range("E6").Value = "  +0.71%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +1.00%  "

$ws.Range("E9").Value = "  -2.19%  "

$ws.Range("E10").Value = "  -1.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.408"
$ws.Range("E11").Value = "  -3.56%  "

$ws.Range("D12").Value = "3.839.75"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.37"
$ws.Range("E14").Value = "  -3.99%  "

$ws.Range("D15").Value = "67.873.88"
$ws.Range("E15").Value = "  -0.99%  "

$ws.Range("E16").Value = "  -1.96%  "

$ws.Range("D17").Value = "3.296.97"
$ws.Range("E17").Value = "  +2.69%  "

$ws.Range("E18").Value = "  -2.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.38"
$ws.Range("E19").Value = "  -1.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "402.97"
$ws.Range("E20").Value = "  +2.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.53"
$ws.Range("E21").Value = "  -2.12%  "

$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.07"
$ws.Range("E23").Value = "  -1.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.508"
$ws.Range("E24").Value = "  -1.41%  "

$ws.Range("E25").Value = "  -1.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.187"
$ws.Range("E26").Value = "  -0.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.46"
$ws.Range("E27").Value = "  -1.40%  "

$ws.Range("E28").Value = "  +0.23%  "

$ws.Range("E29").Value = "  -1.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.64"
$ws.Range("E30").Value = "  -1.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.46"
$ws.Range("E31").Value = "  -3.96%  "

$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.24"
$ws.Range("E34").Value = "  -3.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "164.32"
$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("E36").Value = "  -3.38%  "

$ws.Range("E37").Value = "  -1.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.09"
$ws.Range("E38").Value = "  +2.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.800"
$ws.Range("E39").Value = "  -3.44%  "

$ws.Range("E40").Value = "  -2.59%  "

$ws.Range("E41").Value = "  -3.20%  "

$ws.Range("D42").Value = "2.675.41"
$ws.Range("E42").Value = "  +2.49%  "

$ws.Range("E43").Value = "  -1.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0677"
$ws.Range("E44").Value = "  -1.40%  "

$ws.Range("E45").Value = "  -2.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "334.92"
$ws.Range("E46").Value = "  -3.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.58"
$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("E48").Value = "  -2.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.29"
$ws.Range("E49").Value = "  -0.33%  "

$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.966"
$ws.Range("E51").Value = "  -1.43%  "
